$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels to align with the new Excel file structure
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"

# Update the active selection on the sheet
$ws.Range("C1").Select()
